$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.211.77"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "1.590.37"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'213.83"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'24.00"
$ws.Range("E8").Value = "  +8.11%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "1.818.49"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "1.590.02"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "'3.74"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "28.281.38"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "'63.07"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "'227.11"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").Value = "0.0₃0709"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'9.31"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'151.76"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'15.17"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "'6.57"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "'3.13"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "1.397.93"
$ws.Range("E34").Value = "  -3.78%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -6.76%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'2.54"
$ws.Range("E39").Value = "  +8.33%  "
$ws.Range("D40").Value = "'0.539"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'0.812"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.64"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.87"
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("D45").Value = "'0.987"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "'64.18"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "1.729.43"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").Value = "'2.14"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'87.48"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0524"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -3.26%  "
